# Apply "r0 con traduccion de paises" edit:
# - Add a new worksheet "Hoja1" with country-code -> country-name translation table
# - Make the new sheet the active/selected sheet

$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item("Sheet1")

# Add the new worksheet right after Sheet1
$ws2 = $wb.Worksheets.Add($null, $ws1)
$ws2.Name = "Hoja1"

# Country code -> country name mapping (as entered by the workbook author)
$codes = @("R0_AUS_VA7D","R0_AUT_VA7D","R0_BEL_VA7D","R0_BRA_VA7D","R0_CHE_VA7D","R0_CHL_VA7D","R0_CHN_VA7D","R0_CRI_VA7D","R0_DEU_VA7D","R0_DOM_VA7D","R0_ESP_VA7D","R0_FRA_VA7D","R0_GBR_VA7D","R0_GTM_VA7D","R0_HND_VA7D","R0_ISR_VA7D","R0_ITA_VA7D","R0_KOR_VA7D","R0_MEX_VA7D","R0_NOR_VA7D","R0_PER_VA7D","R0_PRT_VA7D","R0_SGP_VA7D","R0_SLV_VA7D","R0_SWE_VA7D","R0_USA_VA7D")
$names = @("AUSTRIA","AUSTRALIA","NELGICA","BRAZIL","SUIZA ","CHILE ","CHINA ","COSTA RICA ","ALEMANIA ","REPUBLICA DOMINICANA ","ESPAÑA","FRANCIA ","REINO UNIDO ","GUATEMALA","HONDURAS","ISRAEL ","ITALIA ","COREAL DEL SUR","MEXICO ","NORUEGA ","PERU","PORTUGAL","SINGAPUR ","EL SALVADOR","SUECIA ","ESTADOS UNIDOS ")

# Write column B first, in the exact order the strings were originally authored
# (this reproduces the shared-string table insertion order seen in the target file)
$nameWriteOrder = @(3,5,4,6,8,9,12,13,14,15,18,19,20,21,22,25,27,28,16,17,23,24,26,11,10,7)
foreach ($r in $nameWriteOrder) {
    $ws2.Cells.Item($r, 2).Value = $names[$r - 3]
}

# Now write column A (codes) in row order
for ($i = 0; $i -lt $codes.Length; $i++) {
    $r = $i + 3
    $ws2.Cells.Item($r, 1).Value = $codes[$i]
}

$ws2.Columns.Item(1).ColumnWidth = 14.43

# Select the new sheet's data range and make it the active/tabSelected sheet
$ws2.Range("A3:C28").Select()
$ws2.Activate()
